# Excel COM-interop script: updates currentAveragePrice / Leve profit
# columns (H:N) for specific leve rows across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets, reflecting refreshed market-board data
# from the scheduled Sheets runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 45456410
$ws.Range("I112").Value = 850
$ws.Range("J112").Value = 47620960
$ws.Range("K112").Value = 2550
$ws.Range("L112").Value = 142862880
$ws.Range("M112").Value = -1442
$ws.Range("N112").Value = -142865096

# Row 132
$ws.Range("H132").Value = 326841.47
$ws.Range("I132").Value = 337702.88
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 1013108.64
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -1010578.64
$ws.Range("N132").Value = -8060

# Row 137
$ws.Range("H137").Value = 28572892
$ws.Range("I137").Value = 1318.2307
$ws.Range("K137").Value = 3954.6921
$ws.Range("M137").Value = -1404.6921

# Row 141
$ws.Range("H141").Value = 1185.9828
$ws.Range("I141").Value = 604.1739
$ws.Range("J141").Value = 3416.25
$ws.Range("K141").Value = 1812.5217
$ws.Range("L141").Value = 10248.75
$ws.Range("M141").Value = 3367.4783
$ws.Range("N141").Value = -20608.75


$ws = $wb.Worksheets.Item("ARM")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""

# Row 32
$ws.Range("H32").Value = 2611.4949
$ws.Range("I32").Value = 2361.8901
$ws.Range("J32").Value = 3976
$ws.Range("K32").Value = 2361.8901
$ws.Range("L32").Value = 3976
$ws.Range("M32").Value = -2074.8901
$ws.Range("N32").Value = -4550

# Row 88
$ws.Range("H88").Value = 2648.75
$ws.Range("I88").Value = 2750
$ws.Range("J88").Value = 2345
$ws.Range("K88").Value = 2750
$ws.Range("L88").Value = 2345
$ws.Range("M88").Value = -2344
$ws.Range("N88").Value = -3157

# Row 91
$ws.Range("H91").Value = 2648.75
$ws.Range("I91").Value = 2750
$ws.Range("J91").Value = 2345
$ws.Range("K91").Value = 2750
$ws.Range("L91").Value = 2345
$ws.Range("M91").Value = -1346
$ws.Range("N91").Value = -5153


$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2151.5217
$ws.Range("I86").Value = 1927.8387
$ws.Range("J86").Value = 2613.8
$ws.Range("K86").Value = 1927.8387
$ws.Range("L86").Value = 2613.8
$ws.Range("M86").Value = -804.8387
$ws.Range("N86").Value = -4859.8

# Row 89
$ws.Range("H89").Value = 2151.5217
$ws.Range("I89").Value = 1927.8387
$ws.Range("J89").Value = 2613.8
$ws.Range("K89").Value = 9639.193499999999
$ws.Range("L89").Value = 13069
$ws.Range("M89").Value = -4023.193499999999
$ws.Range("N89").Value = -24301


$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 2000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""

# Row 31
$ws.Range("H31").Value = 2242.0815
$ws.Range("I31").Value = 1060.129
$ws.Range("J31").Value = 4277.6665
$ws.Range("K31").Value = 1060.129
$ws.Range("L31").Value = 4277.6665
$ws.Range("M31").Value = -765.1289999999999
$ws.Range("N31").Value = -4867.6665

# Row 34
$ws.Range("H34").Value = 2242.0815
$ws.Range("I34").Value = 1060.129
$ws.Range("J34").Value = 4277.6665
$ws.Range("K34").Value = 1060.129
$ws.Range("L34").Value = 4277.6665
$ws.Range("M34").Value = -858.1289999999999
$ws.Range("N34").Value = -4681.6665

# Row 58
$ws.Range("H58").Value = 1902.5571
$ws.Range("I58").Value = 885.7273
$ws.Range("J58").Value = 3623.3462
$ws.Range("K58").Value = 885.7273
$ws.Range("L58").Value = 3623.3462
$ws.Range("M58").Value = -682.7273
$ws.Range("N58").Value = -4029.3462

# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""

# Row 132
$ws.Range("H132").Value = 2125.1052
$ws.Range("I132").Value = 2116.516
$ws.Range("J132").Value = 2163.1428
$ws.Range("K132").Value = 6349.548000000001
$ws.Range("L132").Value = 6489.428400000001
$ws.Range("M132").Value = -3819.548000000001
$ws.Range("N132").Value = -11549.4284

# Row 134
$ws.Range("H134").Value = 1722.1034
$ws.Range("I134").Value = 1874.9584
$ws.Range("J134").Value = 988.4
$ws.Range("K134").Value = 5624.8752
$ws.Range("L134").Value = 2965.2
$ws.Range("M134").Value = -3089.8752
$ws.Range("N134").Value = -8035.2

# Row 136
$ws.Range("H136").Value = 1902.5571
$ws.Range("I136").Value = 885.7273
$ws.Range("J136").Value = 3623.3462
$ws.Range("K136").Value = 2657.1819
$ws.Range("L136").Value = 10870.0386
$ws.Range("M136").Value = -107.1819
$ws.Range("N136").Value = -15970.0386


$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 3920.8
$ws.Range("I131").Value = 5130
$ws.Range("J131").Value = 3206.2727
$ws.Range("K131").Value = 15390
$ws.Range("L131").Value = 9618.8181
$ws.Range("M131").Value = -10350
$ws.Range("N131").Value = -19698.8181


$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4957.143
$ws.Range("I70").Value = 5266.6665
$ws.Range("J70").Value = 4725
$ws.Range("K70").Value = 5266.6665
$ws.Range("L70").Value = 4725
$ws.Range("M70").Value = -4996.6665
$ws.Range("N70").Value = -5265

# Row 73
$ws.Range("H73").Value = 4957.143
$ws.Range("I73").Value = 5266.6665
$ws.Range("J73").Value = 4725
$ws.Range("K73").Value = 5266.6665
$ws.Range("L73").Value = 4725
$ws.Range("M73").Value = -4330.6665
$ws.Range("N73").Value = -6597

# Row 126
$ws.Range("I126").Value = 2574
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 7722
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -5252
$ws.Range("N126").Value = -9440


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2198.0588
$ws.Range("I16").Value = 1864
$ws.Range("J16").Value = 3757
$ws.Range("K16").Value = 1864
$ws.Range("L16").Value = 3757
$ws.Range("M16").Value = -1694
$ws.Range("N16").Value = -4097

# Row 132
$ws.Range("H132").Value = 4448.706
$ws.Range("I132").Value = 4847.3228
$ws.Range("J132").Value = 329.66666
$ws.Range("K132").Value = 14541.9684
$ws.Range("L132").Value = 988.9999799999999
$ws.Range("M132").Value = -12011.9684
$ws.Range("N132").Value = -6048.99998

# Row 136
$ws.Range("H136").Value = 1589.8096
$ws.Range("I136").Value = 1038.3448
$ws.Range("K136").Value = 3115.0344
$ws.Range("M136").Value = -565.0344000000005


$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3491.8462
$ws.Range("I132").Value = 3901.0667
$ws.Range("J132").Value = 861.1429000000001
$ws.Range("K132").Value = 11703.2001
$ws.Range("L132").Value = 2583.4287
$ws.Range("M132").Value = -9173.2001
$ws.Range("N132").Value = -7643.4287

# Row 136
$ws.Range("H136").Value = 4555
$ws.Range("I136").Value = 4890.622
$ws.Range("J136").Value = 779.25
$ws.Range("K136").Value = 14671.866
$ws.Range("L136").Value = 2337.75
$ws.Range("M136").Value = -12121.866
$ws.Range("N136").Value = -7437.75
